$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Videogame")

# Insert a new column before E (old E -> F), shifting hyperlinked image urls right.
$ws.Columns("E").Insert()

# Header for the new correlative column
$ws.Range("E1").Value = "correlative"

# Correlative numbers per "group" (reset when the console/category changes)
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 6
$ws.Range("E8").Value = 7
$ws.Range("E9").Value = 8
$ws.Range("E10").Value = 9
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 3

# Match the new column's look to column D (which it sits beside)
$ws.Range("E1:E13").Style = $ws.Range("D1:D13").Style
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# Old column E (image urls) is now F; restore its original width behavior.
$ws.Columns("F").ColumnWidth = 90

# Reset the view (no frozen/scrolled topLeftCell, default selection at A1)
$ws.Application.Goto($ws.Range("A1"))
